$wb = $excel.ActiveWorkbook

# Update the "Metadata" sheet (Property/Value pairs in column A/B).
$wsMeta = $wb.Worksheets.Item("Metadata")

$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/coverage-insurance-plan"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# Update the "Elements" sheet.
$wsElements = $wb.Worksheets.Item("Elements")

# AI2 held a duplicate copy of the ele-1/ext-1 constraint text; it is cleared.
$wsElements.Range("AI2").Value = ""

# Q5 mirrors the StructureDefinition URL (same value as Metadata!B2).
$wsElements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/coverage-insurance-plan"
